# "problem135,145" — mark problem 135 (candy) as done and log problem 145 (tree)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

# Row 16 already has problem "135 candy" logged on 2019-03-24 (row for Ray/帅哥Ray).
# Mark it as done.
$ws.Range("E16").Value = "done"

# Add a new row 17 for problem "145 tree", same date as row 16, also done.
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Range("B17").Value = "145 tree"
$ws.Range("E17").Value = "done"

# Leave the cursor where the author left it when saving.
$ws.Range("I11").Select()
